# "added 4wk low sales check" -------------------------------------------
# A 4-week low-sales check was added to the forecasting pipeline. This
# shifted the rolling "Inventory Coverage" (H) down across the board,
# recomputed the "Seasonality Index" (L) values, flagged the final
# forecast week's Reorder Urgency as "Urgent", and bumped the Summary
# sheet's forecast totals.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Forecast Comparison"
$ws2 = $wb.Worksheets.Item(2)   # "Summary"

# --- Forecast Comparison: Inventory Coverage (H) & Seasonality Index (L) ---
$rows = @(
    @{ Row = 2;  H = 15.62;              L = 0.99 },
    @{ Row = 3;  H = 14.62;              L = 1.17 },
    @{ Row = 4;  H = 13.62;              L = 0.92 },
    @{ Row = 5;  H = 12.62;              L = 0.93 },
    @{ Row = 6;  H = 11.62;              L = 1.19 },
    @{ Row = 7;  H = 10.62;              L = 0.9  },
    @{ Row = 8;  H = 9.62;               L = 1.13 },
    @{ Row = 9;  H = 8.62;               L = 0.8  },
    @{ Row = 10; H = 7.62;               L = 0.87 },
    @{ Row = 11; H = 6.62;               L = 1    },
    @{ Row = 12; H = 5.62;               L = 0.96 },
    @{ Row = 13; H = 4.62;               L = 1    },
    @{ Row = 14; H = 3.62;               L = 0.93 },
    @{ Row = 15; H = 2.62;               L = 0.81 },
    @{ Row = 16; H = 1.62;               L = 1.15 },
    @{ Row = 17; H = 0.62;               L = 0.85 }
)

foreach ($r in $rows) {
    $ws1.Cells.Item($r.Row, 8).Value  = $r.H    # column H - Inventory Coverage
    $ws1.Cells.Item($r.Row, 12).Value = $r.L    # column L - Seasonality Index
}

# Row 17 (W25) now needs reordering urgently given the low coverage left.
$ws1.Range("J17").Value = "Urgent"

# --- Summary: forecast totals bumped up by the new check -------------------
$summaryUpdates = @(
    @{ Row = 9;  Value = "26" },
    @{ Row = 10; Value = "13" },
    @{ Row = 11; Value = "6"  },
    @{ Row = 12; Value = "2"  },
    @{ Row = 14; Value = "2"  }
)

foreach ($u in $summaryUpdates) {
    $cell = $ws2.Cells.Item($u.Row, 2)
    # Force text storage (these "Value" column entries are plain text, not
    # numbers, in the source data) and then drop the format stamp so the
    # cell style stays the default, just like its neighbours.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
